# Auto-generated COM-interop script implementing the v1 -> v2 D3FEND library conversion
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Sheet 1: library_content -> library_meta (restructured metadata, dropped 'library_' prefix) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "library_meta"

# Clear existing content first so old rows 11/12 (tab / reference_control_base_urn) are gone
$ws1.Cells.Clear()

$ws1.Cells.Item(1,1).Value = 'type'
$ws1.Cells.Item(1,2).Value = 'library'
$ws1.Cells.Item(2,1).Value = 'urn'
$ws1.Cells.Item(2,2).Value = 'urn:intuitem:risk:library:mitre-d3fend'
$ws1.Cells.Item(3,1).Value = 'version'
$ws1.Cells.Item(3,2).NumberFormat = "@"
$ws1.Cells.Item(3,2).Value = '1'
$ws1.Cells.Item(4,1).Value = 'locale'
$ws1.Cells.Item(4,2).Value = 'en'
$ws1.Cells.Item(5,1).Value = 'publication_date'
$ws1.Cells.Item(5,2).NumberFormat = "@"
$ws1.Cells.Item(5,2).Value = '2025-01-22'
$ws1.Cells.Item(6,1).Value = 'ref_id'
$ws1.Cells.Item(6,2).Value = 'd3fend'
$ws1.Cells.Item(7,1).Value = 'name'
$ws1.Cells.Item(7,2).Value = 'Mitre D3FEND'
$ws1.Cells.Item(8,1).Value = 'description'
$ws1.Cells.Item(8,2).Value = 'A cybersecurity ontology designed to standardize vocabulary for employing techniques to counter malicious cyber threats.
Version - 1.0.0 - 2024-12-20
https://d3fend.mitre.org/resources/'
$ws1.Cells.Item(9,1).Value = 'copyright'
$ws1.Cells.Item(9,2).Value = 'Terms of Use
LICENSE
The MITRE Corporation (MITRE) hereby grants you a non-exclusive, royalty-free license to use D3FEND for research, development, and commercial purposes. Any copy you make for such purposes is authorized provided that you reproduce MITRE’s copyright designation and this license in any such copy.
DISCLAIMERS
ALL DOCUMENTS AND THE INFORMATION CONTAINED THEREIN ARE PROVIDED ON AN "AS IS" BASIS AND THE CONTRIBUTOR, THE ORGANIZATION HE/SHE REPRESENTS OR IS SPONSORED BY (IF ANY), THE MITRE CORPORATION, ITS BOARD OF TRUSTEES, OFFICERS, AGENTS, AND EMPLOYEES, DISCLAIM ALL WARRANTIES, EXPRESS OR IMPLIED, INCLUDING BUT NOT LIMITED TO ANY WARRANTY THAT THE USE OF THE INFORMATION THEREIN WILL NOT INFRINGE ANY RIGHTS OR ANY IMPLIED WARRANTIES OF MERCHANTABILITY OR FITNESS FOR A PARTICULAR PURPOSE.'
$ws1.Cells.Item(10,1).Value = 'provider'
$ws1.Cells.Item(10,2).Value = 'Mitre D3FEND'
$ws1.Cells.Item(11,1).Value = 'packager'
$ws1.Cells.Item(11,2).Value = 'intuitem'

# --- Sheet 2: controls -> controls_meta (small reference-controls metadata sheet) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "controls_meta"

# The old sheet2 held the big reference-controls content table; move it to a new sheet3
# ("controls_content") before overwriting sheet2 with the small meta block.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "controls_content"

# Copy the reference-controls table (ref_id/name/description/category/csf_function) from
# the old controls sheet into the new controls_content sheet, then clear it from ws2.
$srcRange = $ws2.UsedRange
$destRange = $ws3.Range("A1")
$srcRange.Copy($destRange)

$ws2.Cells.Clear()
$ws2.Cells.Item(1,1).Value = 'type'
$ws2.Cells.Item(1,2).Value = 'reference_controls'
$ws2.Cells.Item(2,1).Value = 'base_urn'
$ws2.Cells.Item(2,2).Value = 'urn:intuitem:risk:reference-controls:mitre-d3fend'

$ws1.Range("A1").Select()
